$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All rows 2 through 252 in column C get updated to the new fitness value 7569
$ws.Range("C2:C252").Value = 7569
